$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Create row 78 (new last row). Column A uses the same bold/
#    bordered/centered style as every other row-number cell, so copy
#    that formatting down from A77 before writing the new value.
# ------------------------------------------------------------------
$ws.Cells.Item(77,1).Copy()
$ws.Cells.Item(78,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(78,1).Value = 76
$ws.Cells.Item(78,5).Value = "Cholestasis"

# ------------------------------------------------------------------
# 2) Shift the gene rows 69-77 down into 70-78 (columns B/C/D only;
#    column A and E already hold the correct values for every row).
#    Go from the bottom up so nothing is overwritten before it is read.
# ------------------------------------------------------------------
$ws.Cells.Item(78,2).Value = $ws.Cells.Item(77,2).Value2
$ws.Cells.Item(78,3).Value = $ws.Cells.Item(77,3).Value2
$ws.Cells.Item(78,4).NumberFormat = "@"
$ws.Cells.Item(78,4).Value = $ws.Cells.Item(77,4).Value2
$ws.Cells.Item(77,2).Value = $ws.Cells.Item(76,2).Value2
$ws.Cells.Item(77,3).Value = $ws.Cells.Item(76,3).Value2
$ws.Cells.Item(77,4).NumberFormat = "@"
$ws.Cells.Item(77,4).Value = $ws.Cells.Item(76,4).Value2
$ws.Cells.Item(76,2).Value = $ws.Cells.Item(75,2).Value2
$ws.Cells.Item(76,3).Value = $ws.Cells.Item(75,3).Value2
$ws.Cells.Item(76,4).NumberFormat = "@"
$ws.Cells.Item(76,4).Value = $ws.Cells.Item(75,4).Value2
$ws.Cells.Item(75,2).Value = $ws.Cells.Item(74,2).Value2
$ws.Cells.Item(75,3).Value = $ws.Cells.Item(74,3).Value2
$ws.Cells.Item(75,4).NumberFormat = "@"
$ws.Cells.Item(75,4).Value = $ws.Cells.Item(74,4).Value2
$ws.Cells.Item(74,2).Value = $ws.Cells.Item(73,2).Value2
$ws.Cells.Item(74,3).Value = $ws.Cells.Item(73,3).Value2
$ws.Cells.Item(74,4).NumberFormat = "@"
$ws.Cells.Item(74,4).Value = $ws.Cells.Item(73,4).Value2
$ws.Cells.Item(73,2).Value = $ws.Cells.Item(72,2).Value2
$ws.Cells.Item(73,3).Value = $ws.Cells.Item(72,3).Value2
$ws.Cells.Item(73,4).NumberFormat = "@"
$ws.Cells.Item(73,4).Value = $ws.Cells.Item(72,4).Value2
$ws.Cells.Item(72,2).Value = $ws.Cells.Item(71,2).Value2
$ws.Cells.Item(72,3).Value = $ws.Cells.Item(71,3).Value2
$ws.Cells.Item(72,4).NumberFormat = "@"
$ws.Cells.Item(72,4).Value = $ws.Cells.Item(71,4).Value2
$ws.Cells.Item(71,2).Value = $ws.Cells.Item(70,2).Value2
$ws.Cells.Item(71,3).Value = $ws.Cells.Item(70,3).Value2
$ws.Cells.Item(71,4).NumberFormat = "@"
$ws.Cells.Item(71,4).Value = $ws.Cells.Item(70,4).Value2
$ws.Cells.Item(70,2).Value = $ws.Cells.Item(69,2).Value2
$ws.Cells.Item(70,3).Value = $ws.Cells.Item(69,3).Value2
$ws.Cells.Item(70,4).NumberFormat = "@"
$ws.Cells.Item(70,4).Value = $ws.Cells.Item(69,4).Value2

# ------------------------------------------------------------------
# 3) Write the new VPS50 row at row 69.
# ------------------------------------------------------------------
$ws.Cells.Item(69,2).Value = "VPS50"
$ws.Cells.Item(69,3).Value = "VPS50, EARP/GARPII complex subunit"
$ws.Cells.Item(69,4).NumberFormat = "@"
$ws.Cells.Item(69,4).Value = "2"

# ------------------------------------------------------------------
# 4) Add the new "time_taken" column F: header (bold/bordered style
#    copied from E1) plus a per-row timestamp value.
# ------------------------------------------------------------------
$ws.Cells.Item(1,5).Copy()
$ws.Cells.Item(1,6).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(1,6).Value = "time_taken"

$ws.Cells.Item(2,6).Value = "2021-10-05 13:39:04.780353"
$ws.Cells.Item(3,6).Value = "2021-10-05 13:39:04.780364"
$ws.Cells.Item(4,6).Value = "2021-10-05 13:39:04.780368"
$ws.Cells.Item(5,6).Value = "2021-10-05 13:39:04.780370"
$ws.Cells.Item(6,6).Value = "2021-10-05 13:39:04.780373"
$ws.Cells.Item(7,6).Value = "2021-10-05 13:39:04.780376"
$ws.Cells.Item(8,6).Value = "2021-10-05 13:39:04.780378"
$ws.Cells.Item(9,6).Value = "2021-10-05 13:39:04.780381"
$ws.Cells.Item(10,6).Value = "2021-10-05 13:39:04.780384"
$ws.Cells.Item(11,6).Value = "2021-10-05 13:39:04.780386"
$ws.Cells.Item(12,6).Value = "2021-10-05 13:39:04.780389"
$ws.Cells.Item(13,6).Value = "2021-10-05 13:39:04.780391"
$ws.Cells.Item(14,6).Value = "2021-10-05 13:39:04.780393"
$ws.Cells.Item(15,6).Value = "2021-10-05 13:39:04.780396"
$ws.Cells.Item(16,6).Value = "2021-10-05 13:39:04.780398"
$ws.Cells.Item(17,6).Value = "2021-10-05 13:39:04.780401"
$ws.Cells.Item(18,6).Value = "2021-10-05 13:39:04.780404"
$ws.Cells.Item(19,6).Value = "2021-10-05 13:39:04.780406"
$ws.Cells.Item(20,6).Value = "2021-10-05 13:39:04.780409"
$ws.Cells.Item(21,6).Value = "2021-10-05 13:39:04.780411"
$ws.Cells.Item(22,6).Value = "2021-10-05 13:39:04.780414"
$ws.Cells.Item(23,6).Value = "2021-10-05 13:39:04.780416"
$ws.Cells.Item(24,6).Value = "2021-10-05 13:39:04.780419"
$ws.Cells.Item(25,6).Value = "2021-10-05 13:39:04.780421"
$ws.Cells.Item(26,6).Value = "2021-10-05 13:39:04.780424"
$ws.Cells.Item(27,6).Value = "2021-10-05 13:39:04.780427"
$ws.Cells.Item(28,6).Value = "2021-10-05 13:39:04.780429"
$ws.Cells.Item(29,6).Value = "2021-10-05 13:39:04.780432"
$ws.Cells.Item(30,6).Value = "2021-10-05 13:39:04.780435"
$ws.Cells.Item(31,6).Value = "2021-10-05 13:39:04.780437"
$ws.Cells.Item(32,6).Value = "2021-10-05 13:39:04.780440"
$ws.Cells.Item(33,6).Value = "2021-10-05 13:39:04.780442"
$ws.Cells.Item(34,6).Value = "2021-10-05 13:39:04.780445"
$ws.Cells.Item(35,6).Value = "2021-10-05 13:39:04.780448"
$ws.Cells.Item(36,6).Value = "2021-10-05 13:39:04.780450"
$ws.Cells.Item(37,6).Value = "2021-10-05 13:39:04.780453"
$ws.Cells.Item(38,6).Value = "2021-10-05 13:39:04.780456"
$ws.Cells.Item(39,6).Value = "2021-10-05 13:39:04.780458"
$ws.Cells.Item(40,6).Value = "2021-10-05 13:39:04.780461"
$ws.Cells.Item(41,6).Value = "2021-10-05 13:39:04.780463"
$ws.Cells.Item(42,6).Value = "2021-10-05 13:39:04.780466"
$ws.Cells.Item(43,6).Value = "2021-10-05 13:39:04.780468"
$ws.Cells.Item(44,6).Value = "2021-10-05 13:39:04.780471"
$ws.Cells.Item(45,6).Value = "2021-10-05 13:39:04.780473"
$ws.Cells.Item(46,6).Value = "2021-10-05 13:39:04.780476"
$ws.Cells.Item(47,6).Value = "2021-10-05 13:39:04.780478"
$ws.Cells.Item(48,6).Value = "2021-10-05 13:39:04.780481"
$ws.Cells.Item(49,6).Value = "2021-10-05 13:39:04.780483"
$ws.Cells.Item(50,6).Value = "2021-10-05 13:39:04.780486"
$ws.Cells.Item(51,6).Value = "2021-10-05 13:39:04.780488"
$ws.Cells.Item(52,6).Value = "2021-10-05 13:39:04.780491"
$ws.Cells.Item(53,6).Value = "2021-10-05 13:39:04.780493"
$ws.Cells.Item(54,6).Value = "2021-10-05 13:39:04.780496"
$ws.Cells.Item(55,6).Value = "2021-10-05 13:39:04.780499"
$ws.Cells.Item(56,6).Value = "2021-10-05 13:39:04.780501"
$ws.Cells.Item(57,6).Value = "2021-10-05 13:39:04.780504"
$ws.Cells.Item(58,6).Value = "2021-10-05 13:39:04.780506"
$ws.Cells.Item(59,6).Value = "2021-10-05 13:39:04.780509"
$ws.Cells.Item(60,6).Value = "2021-10-05 13:39:04.780511"
$ws.Cells.Item(61,6).Value = "2021-10-05 13:39:04.780514"
$ws.Cells.Item(62,6).Value = "2021-10-05 13:39:04.780516"
$ws.Cells.Item(63,6).Value = "2021-10-05 13:39:04.780519"
$ws.Cells.Item(64,6).Value = "2021-10-05 13:39:04.780521"
$ws.Cells.Item(65,6).Value = "2021-10-05 13:39:04.780524"
$ws.Cells.Item(66,6).Value = "2021-10-05 13:39:04.780527"
$ws.Cells.Item(67,6).Value = "2021-10-05 13:39:04.780530"
$ws.Cells.Item(68,6).Value = "2021-10-05 13:39:04.780532"
$ws.Cells.Item(69,6).Value = "2021-10-05 13:39:04.780535"
$ws.Cells.Item(70,6).Value = "2021-10-05 13:39:04.780538"
$ws.Cells.Item(71,6).Value = "2021-10-05 13:39:04.780540"
$ws.Cells.Item(72,6).Value = "2021-10-05 13:39:04.780543"
$ws.Cells.Item(73,6).Value = "2021-10-05 13:39:04.780545"
$ws.Cells.Item(74,6).Value = "2021-10-05 13:39:04.780548"
$ws.Cells.Item(75,6).Value = "2021-10-05 13:39:04.780550"
$ws.Cells.Item(76,6).Value = "2021-10-05 13:39:04.780553"
$ws.Cells.Item(77,6).Value = "2021-10-05 13:39:04.780555"
$ws.Cells.Item(78,6).Value = "2021-10-05 13:39:04.780559"
